# Rspo3-Fzd8.xlsx — "Natmi following Dr Hou advice"
#
# Recomputed NATMI ligand-receptor edge-weight statistics for the
# FAPs(Rspo3) -> Fzd8 pair and extended the Target-cluster coverage
# from 3 clusters (ECs, FAPs, sCs) to 5 clusters, adding M1 and
# Neutro as new target-cluster rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster / Ligand symbol / Receptor symbol are identical
# (FAPs / Rspo3 / Fzd8) for every data row; only the Target cluster
# (col D) and the computed statistics (cols E:T) vary per row.
$targetClusters = @("ECs", "FAPs", "M1", "Neutro", "sCs")
for ($i = 0; $i -lt $targetClusters.Length; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = "FAPs"
    $ws.Range("B$r").Value = "Rspo3"
    $ws.Range("C$r").Value = "Fzd8"
    $ws.Range("D$r").Value = $targetClusters[$i]
}

# Row 2: Target cluster = ECs
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.822099333333334
$ws.Range("H2").Value = 8.466298
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.416422666666667
$ws.Range("N2").Value = 7.249268
$ws.Range("O2").Value = 0.2729115228630338
$ws.Range("P2").Value = 0.2729115228630338
$ws.Range("Q2").Value = 6.819384796651556
$ws.Range("R2").Value = 61.374463169864
$ws.Range("S2").Value = 0.2729115228630338
$ws.Range("T2").Value = 0.2729115228630338

# Row 3: Target cluster = FAPs
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.822099333333334
$ws.Range("H3").Value = 8.466298
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.913654666666667
$ws.Range("N3").Value = 8.740964
$ws.Range("O3").Value = 0.3290690586319826
$ws.Range("P3").Value = 0.3290690586319827
$ws.Range("Q3").Value = 8.222622892363555
$ws.Range("R3").Value = 74.003606031272
$ws.Range("S3").Value = 0.3290690586319826
$ws.Range("T3").Value = 0.3290690586319827

# Row 4: Target cluster = M1
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.822099333333334
$ws.Range("H4").Value = 8.466298
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.02921733333333333
$ws.Range("N4").Value = 0.08765200000000001
$ws.Range("O4").Value = 0.003299814657423431
$ws.Range("P4").Value = 0.003299814657423432
$ws.Range("Q4").Value = 0.08245421692177779
$ws.Range("R4").Value = 0.7420879522960001
$ws.Range("S4").Value = 0.003299814657423431
$ws.Range("T4").Value = 0.003299814657423432

# Row 5: Target cluster = Neutro
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.822099333333334
$ws.Range("H5").Value = 8.466298
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.02851766666666667
$ws.Range("N5").Value = 0.085553
$ws.Range("O5").Value = 0.003220794087830817
$ws.Range("P5").Value = 0.003220794087830818
$ws.Range("Q5").Value = 0.08047968808822223
$ws.Range("R5").Value = 0.7243171927940001
$ws.Range("S5").Value = 0.003220794087830817
$ws.Range("T5").Value = 0.003220794087830818

# Row 6: Target cluster = sCs
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.822099333333334
$ws.Range("H6").Value = 8.466298
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.466422333333333
$ws.Range("N6").Value = 10.399267
$ws.Range("O6").Value = 0.3914988097597293
$ws.Range("P6").Value = 0.3914988097597293
$ws.Range("Q6").Value = 9.782588155951778
$ws.Range("R6").Value = 88.043293403566
$ws.Range("S6").Value = 0.3914988097597293
$ws.Range("T6").Value = 0.3914988097597293
